# Saldo.xlsx update: remove a handful of rows that dropped out of the
# export and move the GILSON (004474776) account up into the big-balance
# section with its refreshed balance.
#
# Sheet layout: row 1 = header (Conta / Nome / Saldo), data starts at row 2.
# Row numbers below are the *original* positions in the sheet before any
# deletes happen in this script.
#   row 5  -> 004241147 ANTONIO    64312.73   (remove entirely)
#   row 7  -> 004415557 FILIPE     37886.67   (remove entirely)
#   row 8  -> 004983395 MARCELO    37831.02   (remove entirely)
#   row 18 -> 005219257 CAROLINE   1029.25    (remove entirely)
#   row 35 -> 004474776 GILSON     202.65     (remove; re-inserted above with
#                                              new balance 24202.65)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()

# After the five deletes above, THIAGO (005064129) sits on row 6 and
# PHYLIA (004690692) on row 7. Insert a new row 7 for GILSON's refreshed
# balance so it lands between them.
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "004474776"
$ws.Cells.Item(7, 1).ClearFormats()
$ws.Cells.Item(7, 2).Value = "GILSON"
$ws.Cells.Item(7, 3).Value = 24202.65
